# Fix the "tile" column (column B) in the Potluck move/tile/pass-go table
# to hold the board position number instead of the tile's name, for the
# three data rows (19-21) under the headers in row 18.
#
#   B19: "old creek" -> 2
#   B20: "go"         -> 1
#   B21: "jail"        -> 11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B19").Value = 2
$ws.Range("B20").Value = 1
$ws.Range("B21").Value = 11
